$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the sheet so we can rebuild it with the new "keyword" column
# and the additional rows described in the commit.
$ws.UsedRange.ClearContents()

# --- Cells whose text re-uses strings already present in the workbook ---
$ws.Range("B1").Value = "prefix"
$ws.Range("C1").Value = "suffix"
$ws.Range("D1").Value = "parent_prefix"
$ws.Range("E1").Value = "context"
$ws.Range("F1").Value = "group"
$ws.Range("G1").Value = "type"
$ws.Range("H1").Value = "subtype"
$ws.Range("I1").Value = "comments"
$ws.Range("B2").Value = "dat"
$ws.Range("E2").Value = "PostgreSQL"
$ws.Range("G2").Value = "table"
$ws.Range("I2").Value = "Network Object (w/ identifiers) data tables"
$ws.Range("B3").Value = "geom"
$ws.Range("E3").Value = "PostgreSQL"
$ws.Range("G3").Value = "table"
$ws.Range("I3").Value = "Network Object PostGIS geometry tables"
$ws.Range("B4").Value = "mtx"
$ws.Range("E4").Value = "PostgreSQL"
$ws.Range("G4").Value = "table"
$ws.Range("I4").Value = "Skim/Demand matrices"
$ws.Range("B5").Value = "net"
$ws.Range("E5").Value = "PostgreSQL"
$ws.Range("G5").Value = "table"
$ws.Range("I5").Value = "Network Object identifier tables"
$ws.Range("B6").Value = "gfx"
$ws.Range("E6").Value = "PostgreSQL"
$ws.Range("G6").Value = "table"
$ws.Range("B7").Value = "tim"
$ws.Range("E7").Value = "PostgreSQL"
$ws.Range("G7").Value = "function"
$ws.Range("I7").Value = "General purpose PostgreSQL function"
$ws.Range("B8").Value = "gfx"
$ws.Range("D8").Value = "tim"
$ws.Range("E8").Value = "PostgreSQL"
$ws.Range("G8").Value = "function"
$ws.Range("I8").Value = "PostGIS specific PostgreSQL function"
$ws.Range("A9").Value = "b"
$ws.Range("E9").Value = "PHP"
$ws.Range("G9").Value = "GET keyword"
$ws.Range("A10").Value = "d"
$ws.Range("E10").Value = "PHP"
$ws.Range("G10").Value = "GET keyword"
$ws.Range("E11").Value = "PHP"
$ws.Range("F11").Value = "d"
$ws.Range("G11").Value = "GET keyword"
$ws.Range("E12").Value = "PHP"
$ws.Range("F12").Value = "d"
$ws.Range("G12").Value = "GET keyword"
$ws.Range("A13").Value = "g"
$ws.Range("E13").Value = "PHP"
$ws.Range("G13").Value = "GET keyword"
$ws.Range("A14").Value = "gpt"
$ws.Range("E14").Value = "PHP"
$ws.Range("G14").Value = "GET keyword"
$ws.Range("H14").Value = "geometry"
$ws.Range("I14").Value = "Geometry - Point"
$ws.Range("A15").Value = "gln"
$ws.Range("E15").Value = "PHP"
$ws.Range("G15").Value = "GET keyword"
$ws.Range("H15").Value = "geometry"
$ws.Range("I15").Value = "Geometry - Line"
$ws.Range("A16").Value = "gpg"
$ws.Range("E16").Value = "PHP"
$ws.Range("G16").Value = "GET keyword"
$ws.Range("H16").Value = "geometry"
$ws.Range("I16").Value = "Geometry - Polygon"

# --- Cells that introduce the new shared strings for this edit ---
$ws.Range("I6").Value = "Helper Graphics/PostGIS tables"
$ws.Range("A1").Value = "keyword"
$ws.Range("I11").Value = "Count of requested attributes"
$ws.Range("I12").Value = "Requested attributes, e.g. bn=2&b0=hi&b1=there"
$ws.Range("A11").Value = "dn"
$ws.Range("A12").Value = "d[0-9]+"
$ws.Range("I10").Value = "Data request"
$ws.Range("I9").Value = "Identifier request"
$ws.Range("I13").Value = "Geometry request"

# Restore the selection reported for this sheet after the edit
$ws.Range("I14").Select()

